$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.402671933174133
$ws.Range("B1").Value = 2.417877912521362
$ws.Range("C1").Value = 3.074418783187866
$ws.Range("D1").Value = 3.568416357040405
$ws.Range("E1").Value = 1.870812892913818
